$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 13.43319801795061
$ws.Range("D2").Value = 5.025506521011489
$ws.Range("E2").Value = 13.79870895180907
$ws.Range("F2").Value = 23.73888337064486
$ws.Range("G2").Value = 27.78686553745599
$ws.Range("H2").Value = 14.12461759485033
$ws.Range("K2").Value = 16.65151356881643
$ws.Range("L2").Value = 9.202506850619693
$ws.Range("O2").Value = 21.35583386709204
$ws.Range("C3").Value = 13.33295043205572
$ws.Range("D3").Value = 4.968424097414155
$ws.Range("E3").Value = 13.73667503250987
$ws.Range("F3").Value = 23.8542960704896
$ws.Range("G3").Value = 27.98563544677185
$ws.Range("H3").Value = 14.20283793606021
$ws.Range("K3").Value = 15.91180331386633
$ws.Range("L3").Value = 9.204210684375211
$ws.Range("O3").Value = 21.49593963771159
$ws.Range("C4").Value = 13.27440995492887
$ws.Range("D4").Value = 4.932864739493822
$ws.Range("E4").Value = 13.70166410780223
$ws.Range("F4").Value = 23.93395073440352
$ws.Range("G4").Value = 28.12111684914884
$ws.Range("H4").Value = 14.25400941028919
$ws.Range("K4").Value = 15.43885356262375
$ws.Range("L4").Value = 9.206772773886609
$ws.Range("O4").Value = 21.58851270462081
$ws.Range("C5").Value = 13.2513343797739
$ws.Range("D5").Value = 4.91825512949019
$ws.Range("E5").Value = 13.68818202697163
$ws.Range("F5").Value = 23.96860738416449
$ws.Range("G5").Value = 28.17967415051406
$ws.Range("H5").Value = 14.27565216995571
$ws.Range("K5").Value = 15.24161966241907
$ws.Range("L5").Value = 9.208198509082917
$ws.Range("O5").Value = 21.62787751976444
$ws.Range("C6").Value = 13.24755045382205
$ws.Range("D6").Value = 4.915822314480311
$ws.Range("E6").Value = 13.68599106328844
$ws.Range("F6").Value = 23.97449436096376
$ws.Range("G6").Value = 28.18959881685347
$ws.Range("H6").Value = 14.27929362967479
$ws.Range("K6").Value = 15.2086040793515
$ws.Range("L6").Value = 9.20845831392883
$ws.Range("O6").Value = 21.63451293146479
$ws.Range("C7").Value = 13.2740955618652
$ws.Range("D7").Value = 4.932668177917066
$ws.Range("E7").Value = 13.70147909077273
$ws.Range("F7").Value = 23.93440925017522
$ws.Range("G7").Value = 28.12189306061071
$ws.Range("H7").Value = 14.25429809432054
$ws.Range("K7").Value = 15.43621152171547
$ws.Range("L7").Value = 9.20679045599061
$ws.Range("O7").Value = 21.58903695755878
$ws.Range("C8").Value = 13.39802161099211
$ws.Range("D8").Value = 5.005936108467584
$ws.Range("E8").Value = 13.77668778292745
$ws.Range("F8").Value = 23.77684496305857
$ws.Range("G8").Value = 27.85259420519876
$ws.Range("H8").Value = 14.1509351009285
$ws.Range("K8").Value = 16.40047573125465
$ws.Range("L8").Value = 9.202780021005891
$ws.Range("O8").Value = 21.40277982016207
$ws.Range("C9").Value = 13.66387255244034
$ws.Range("D9").Value = 5.145155475848163
$ws.Range("E9").Value = 13.94805821406249
$ws.Range("F9").Value = 23.53821091804529
$ws.Range("G9").Value = 27.4324800205476
$ws.Range("H9").Value = 13.97321226305442
$ws.Range("K9").Value = 18.13455813549566
$ws.Range("L9").Value = 9.206918663590709
$ws.Range("O9").Value = 21.08974910675762
$ws.Range("C10").Value = 13.87160136469021
$ws.Range("D10").Value = 5.244189697498191
$ws.Range("E10").Value = 14.08773986681142
$ws.Range("F10").Value = 23.4065613298959
$ws.Range("G10").Value = 27.19149939364818
$ws.Range("H10").Value = 13.85789980352174
$ws.Range("K10").Value = 19.30419622444174
$ws.Range("L10").Value = 9.217236395642265
$ws.Range("O10").Value = 20.89197036954377
$ws.Range("C11").Value = 13.96845510094231
$ws.Range("D11").Value = 5.288428713771115
$ws.Range("E11").Value = 14.15408949993576
$ws.Range("F11").Value = 23.3563123180338
$ws.Range("G11").Value = 27.09695249776254
$ws.Range("H11").Value = 13.80876390007036
$ws.Range("K11").Value = 19.812228573458
$ws.Range("L11").Value = 9.223498870974957
$ws.Range("O11").Value = 20.80907083916158
$ws.Range("C12").Value = 14.00544025873686
$ws.Range("D12").Value = 5.305055304432815
$ws.Range("E12").Value = 14.17960117576173
$ws.Range("F12").Value = 23.33868253038955
$ws.Range("G12").Value = 27.06334886553498
$ws.Range("H12").Value = 13.79063581839381
$ws.Range("K12").Value = 20.00105810024045
$ws.Range("L12").Value = 9.226094662543813
$ws.Range("O12").Value = 20.77870285737411
$ws.Range("C13").Value = 13.99746159652836
$ws.Range("D13").Value = 5.301480203885816
$ws.Range("E13").Value = 14.17408987170619
$ws.Range("F13").Value = 23.34241703596632
$ws.Range("G13").Value = 27.07048770892707
$ws.Range("H13").Value = 13.7945187193176
$ws.Range("K13").Value = 19.96054955811304
$ws.Range("L13").Value = 9.225525657254749
$ws.Range("O13").Value = 20.78519747061768
$ws.Range("C14").Value = 13.97149187869805
$ws.Range("D14").Value = 5.289799164285375
$ws.Range("E14").Value = 14.1561807158648
$ws.Range("F14").Value = 23.35483381000205
$ws.Range("G14").Value = 27.09414365485092
$ws.Range("H14").Value = 13.80726289278914
$ws.Range("K14").Value = 19.8278353480906
$ws.Range("L14").Value = 9.223707940285008
$ws.Range("O14").Value = 20.80655187588351
$ws.Range("C15").Value = 13.95562398198389
$ws.Range("D15").Value = 5.282627550508355
$ws.Range("E15").Value = 14.14526065209841
$ws.Range("F15").Value = 23.36262191345931
$ws.Range("G15").Value = 27.10892092256303
$ws.Range("H15").Value = 13.8151314327653
$ws.Range("K15").Value = 19.74607905848043
$ws.Range("L15").Value = 9.22262370985073
$ws.Range("O15").Value = 20.81976567982561
$ws.Range("C16").Value = 13.86531660386784
$ws.Range("D16").Value = 5.24128145707528
$ws.Range("E16").Value = 14.08345880384061
$ws.Range("F16").Value = 23.4100402529971
$ws.Range("G16").Value = 27.19798432156321
$ws.Range("H16").Value = 13.86117786150963
$ws.Range("K16").Value = 19.27050291596231
$ws.Range("L16").Value = 9.216858592543044
$ws.Range("O16").Value = 20.89753097119914
$ws.Range("C17").Value = 13.81049829971005
$ws.Range("D17").Value = 5.215702477384736
$ws.Range("E17").Value = 14.04625246678165
$ws.Range("F17").Value = 23.44160764486516
$ws.Range("G17").Value = 27.25650574859248
$ws.Range("H17").Value = 13.89027707383807
$ws.Range("K17").Value = 18.97252246073845
$ws.Range("L17").Value = 9.213722849518566
$ws.Range("O17").Value = 20.94705378822214
$ws.Range("C18").Value = 13.77919194179887
$ws.Range("D18").Value = 5.200914200828429
$ws.Range("E18").Value = 14.02511772453844
$ws.Range("F18").Value = 23.46067094802112
$ws.Range("G18").Value = 27.29158296425252
$ws.Range("H18").Value = 13.90732662148194
$ws.Range("K18").Value = 18.79887355652386
$ws.Range("L18").Value = 9.212067010098499
$ws.Range("O18").Value = 20.97620310489464
$ws.Range("C19").Value = 13.76863147524469
$ws.Range("D19").Value = 5.195894382804901
$ws.Range("E19").Value = 14.01800795309943
$ws.Range("F19").Value = 23.46728077333363
$ws.Range("G19").Value = 27.3037019536073
$ws.Range("H19").Value = 13.91315293460528
$ws.Range("K19").Value = 18.73969426973193
$ws.Range("L19").Value = 9.211531785142931
$ws.Range("O19").Value = 20.98618656912442
$ws.Range("C20").Value = 13.8163108686402
$ws.Range("D20").Value = 5.218433327242206
$ws.Range("E20").Value = 14.0501858038802
$ws.Range("F20").Value = 23.43815332357505
$ws.Range("G20").Value = 27.25012913487537
$ws.Range("H20").Value = 13.88714706732019
$ws.Range("K20").Value = 19.00447746046661
$ws.Range("L20").Value = 9.214041371183731
$ws.Range("O20").Value = 20.94171311048806
$ws.Range("C21").Value = 13.97911167757828
$ws.Range("D21").Value = 5.293233652080723
$ws.Range("E21").Value = 14.16143072016116
$ws.Range("F21").Value = 23.35114865464108
$ws.Range("G21").Value = 27.08713539873062
$ws.Range("H21").Value = 13.80350662099614
$ws.Range("K21").Value = 19.86691378096278
$ws.Range("L21").Value = 9.224235770537952
$ws.Range("O21").Value = 20.80025171085428
$ws.Range("C22").Value = 14.08729677696336
$ws.Range("D22").Value = 5.341382982297663
$ws.Range("E22").Value = 14.23638010708071
$ws.Range("F22").Value = 23.30244239675661
$ws.Range("G22").Value = 26.99344223113803
$ws.Range("H22").Value = 13.75163309054258
$ws.Range("K22").Value = 20.40983321929972
$ws.Range("L22").Value = 9.232205257406036
$ws.Range("O22").Value = 20.71377229740573
$ws.Range("C23").Value = 14.02940306561503
$ws.Range("D23").Value = 5.315755104031731
$ws.Range("E23").Value = 14.19617871766742
$ws.Range("F23").Value = 23.32768763630002
$ws.Range("G23").Value = 27.04226374061602
$ws.Range("H23").Value = 13.77906325496765
$ws.Range("K23").Value = 20.12199100266949
$ws.Range("L23").Value = 9.227832675166418
$ws.Range("O23").Value = 20.75937888640503
$ws.Range("C24").Value = 13.81368235164775
$ws.Range("D24").Value = 5.217198966545896
$ws.Range("E24").Value = 14.04840674448635
$ws.Range("F24").Value = 23.43971217362472
$ws.Range("G24").Value = 27.25300754220785
$ws.Range("H24").Value = 13.88856114646797
$ws.Range("K24").Value = 18.99003787384024
$ws.Range("L24").Value = 9.213896909738546
$ws.Range("O24").Value = 20.94412551950426
$ws.Range("C25").Value = 13.5896603214303
$ws.Range("D25").Value = 5.108026032393994
$ws.Range("E25").Value = 13.89921814801159
$ws.Range("F25").Value = 23.59515308531045
$ws.Range("G25").Value = 27.53437932032314
$ws.Range("H25").Value = 14.01861398456477
$ws.Range("K25").Value = 17.6832330814057
$ws.Range("L25").Value = 9.204517623164008
$ws.Range("O25").Value = 21.16880399836092
